$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Förändrad) for rows 2-14 changed from 45184 to 45185 (date +1 day)
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value2 = 45185
    }
}
